# Jogos_da_Semana_FlashScore_2024-10-31.xlsx
# A new match (QATAR - QSL: Al-Sadd vs Al Rayyan) was inserted as row 5,
# pushing the two following Saudi Professional League rows down by one.
# Several odds columns on the EGYPT - PREMIER LEAGUE row (row 2) and on the
# two shifted rows were refreshed with updated market prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 5 (shifts old rows 5->6, 6->7)
$ws.Rows.Item(5).Insert()

# 2) Update odds in row 2 (EGYPT - PREMIER LEAGUE: El Gouna vs ZED)
$ws.Range("G2").Value = 3.6
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 2.27
$ws.Range("J2").Value = 4.1
$ws.Range("K2").Value = 1.95
$ws.Range("L2").Value = 2.82
$ws.Range("N2").Value = 5.5
$ws.Range("O2").Value = 1.45
$ws.Range("P2").Value = 2.55
$ws.Range("R2").Value = 1.55
$ws.Range("S2").Value = 1.47
$ws.Range("T2").Value = 2.52
$ws.Range("U2").Value = 1.9
$ws.Range("V2").Value = 1.8
$ws.Range("W2").Value = 8.5
$ws.Range("X2").Value = 18.5
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 55
$ws.Range("AA2").Value = 37
$ws.Range("AB2").Value = 45
$ws.Range("AC2").Value = 5.5
$ws.Range("AD2").Value = 5.4
$ws.Range("AE2").Value = 14.5
$ws.Range("AF2").Value = 80
$ws.Range("AH2").Value = 6.3
$ws.Range("AI2").Value = 10.25
$ws.Range("AJ2").Value = 9
$ws.Range("AK2").Value = 24
$ws.Range("AL2").Value = 21
$ws.Range("AM2").Value = 35
$ws.Range("AN2").Value = 5.4
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 26
$ws.Range("AQ2").Value = 110
$ws.Range("AR2").Value = 150
$ws.Range("AS2").Value = 300
$ws.Range("AT2").Value = 2.52
$ws.Range("AU2").Value = 6.6
$ws.Range("AW2").Value = 4.1
$ws.Range("AX2").Value = 12
$ws.Range("AY2").Value = 19.5
$ws.Range("AZ2").Value = 50
$ws.Range("BA2").Value = 80

# 3) Populate the newly inserted row 5 (QATAR - QSL: Al-Sadd vs Al Rayyan)
$ws.Range("A5").Value = "K6I8gpX9"
$ws.Range("B5").Value = "31/10/2024"
$ws.Range("C5").Value = "13:30"
$ws.Range("D5").Value = "QATAR - QSL"
$ws.Range("E5").Value = "Al-Sadd"
$ws.Range("F5").Value = "Al Rayyan"
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 4.05
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 2.4
$ws.Range("K5").Value = 2.6
$ws.Range("L5").Value = 3.25
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.08
$ws.Range("P5").Value = 6.3
$ws.Range("Q5").Value = 1.28
$ws.Range("R5").Value = 3.35
$ws.Range("S5").Value = 1.18
$ws.Range("T5").Value = 4.25
$ws.Range("U5").Value = 1.29
$ws.Range("V5").Value = 3.25
$ws.Range("W5").Value = 18.5
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 24
$ws.Range("AA5").Value = 13.5
$ws.Range("AB5").Value = 14.5
$ws.Range("AC5").Value = 10.75
$ws.Range("AD5").Value = 10.25
$ws.Range("AE5").Value = 10.5
$ws.Range("AF5").Value = 23
$ws.Range("AG5").Value = 90
$ws.Range("AH5").Value = 22
$ws.Range("AI5").Value = 26
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 45
$ws.Range("AL5").Value = 21
$ws.Range("AM5").Value = 18.5
$ws.Range("AN5").Value = 5
$ws.Range("AO5").Value = 9.75
$ws.Range("AP5").Value = 11.75
$ws.Range("AQ5").Value = 27
$ws.Range("AR5").Value = 32
$ws.Range("AS5").Value = 75
$ws.Range("AT5").Value = 4.25
$ws.Range("AU5").Value = 5.9
$ws.Range("AV5").Value = 27
$ws.Range("AW5").Value = 6.1
$ws.Range("AX5").Value = 15
$ws.Range("AY5").Value = 14.5
$ws.Range("AZ5").Value = 50
$ws.Range("BA5").Value = 50
$ws.Range("BB5").Value = 100
$ws.Range("BC5").Value = 250
$ws.Range("BD5").Value = 51

# 4) Minor odds tweak carried with the shifted row 6 (Al Shabab vs Al Wehda)
$ws.Range("AG6").Value = 700

# 5) Minor odds tweaks carried with the shifted row 7 (Al Okhdood vs Al Orubah)
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 8.5
$ws.Range("Q7").Value = 1.83
$ws.Range("R7").Value = 1.98
